$wb = $excel.ActiveWorkbook

# "Generate Report for Handback" - refresh the handoff/handback timestamps
# for the first data row (5fe23c8a-...) on the zh-cn and de-de sheets.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-03 06:51:35"
$wsZhCn.Range("K2").Value = "2016-09-03 06:51:52"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-03 06:51:40"
$wsDeDe.Range("K2").Value = "2016-09-03 06:51:59"
